$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting rows 39-46 down to 40-47.
$ws.Rows.Item(39).Insert()

# Match the number format of the date cell (same as the rest of column D)
$ws.Range("D39").NumberFormat = $ws.Range("D40").NumberFormat

# Fill in the new row 39 data (a new weekly price record for the same market/origin)
$ws.Cells.Item(39, 1).Value = 7
$ws.Cells.Item(39, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(39, 3).Value = "Ñuble"
$ws.Cells.Item(39, 4).Value = 44559
$ws.Cells.Item(39, 5).Value = 16
$ws.Cells.Item(39, 6).Value = 100112022
$ws.Cells.Item(39, 7).Value = "Arveja Verde"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 60
$ws.Cells.Item(39, 11).Value = 18000
$ws.Cells.Item(39, 12).Value = 19000
$ws.Cells.Item(39, 13).Value = 18500
$ws.Cells.Item(39, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(39, 16).Value = 740
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
